# Update the generated-report timestamps on the handback-status workbook.
# The "zh-cn" and "de-de" sheets each have a row (row 3, for the
# 41711229-... handoff file) whose "Correspond Handoff Datetime" (column D)
# and "Correspond Handback DateTime" (column G) values need to be refreshed
# to reflect a newer report generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-20 03:30:09"
$wsZhCn.Range("G3").Value = "2016-01-20 03:30:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-20 03:30:20"
$wsDeDe.Range("G3").Value = "2016-01-20 03:31:19"
